$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Sat Feb 17 22:26:44 EST 2024"
$ws.Range("B3").Value = "Sat Feb 17 22:26:54 EST 2024"
$ws.Range("B4").Value = "Sat Feb 17 22:27:03 EST 2024"
$ws.Range("B5").Value = "Sat Feb 17 22:27:13 EST 2024"
$ws.Range("B6").Value = "Sat Feb 17 22:27:23 EST 2024"
$ws.Range("B7").Value = "Sat Feb 17 22:27:33 EST 2024"
$ws.Range("B8").Value = "Sat Feb 17 22:27:42 EST 2024"
$ws.Range("B9").Value = "Sat Feb 17 22:27:52 EST 2024"
$ws.Range("B10").Value = "Sat Feb 17 22:28:02 EST 2024"
$ws.Range("B13").Value = "Sat Feb 17 22:28:12 EST 2024"
$ws.Range("B14").Value = "Sat Feb 17 22:28:22 EST 2024"
$ws.Range("B15").Value = "Sat Feb 17 22:28:31 EST 2024"
$ws.Range("B16").Value = "Sat Feb 17 22:28:41 EST 2024"
$ws.Range("B17").Value = "Sat Feb 17 22:28:51 EST 2024"
$ws.Range("B18").Value = "Sat Feb 17 22:29:01 EST 2024"
$ws.Range("B19").Value = "Sat Feb 17 22:29:10 EST 2024"
$ws.Range("B20").Value = "Sat Feb 17 22:29:20 EST 2024"
$ws.Range("B21").Value = "Sat Feb 17 22:29:30 EST 2024"
$ws.Range("B22").Value = "Sat Feb 17 22:29:40 EST 2024"
$ws.Range("B23").Value = "Sat Feb 17 22:29:50 EST 2024"
$ws.Range("B24").Value = "Sat Feb 17 22:30:00 EST 2024"
$ws.Range("B25").Value = "Sat Feb 17 22:30:10 EST 2024"
$ws.Range("B26").Value = "Sat Feb 17 22:30:19 EST 2024"
$ws.Range("B27").Value = "Sat Feb 17 22:30:29 EST 2024"
$ws.Range("B28").Value = "Sat Feb 17 22:30:39 EST 2024"
$ws.Range("B29").Value = "Sat Feb 17 22:30:49 EST 2024"
$ws.Range("B30").Value = "Sat Feb 17 22:30:58 EST 2024"
